$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new transaction row is being inserted at row 2, pushing the existing
# rows (2-5) down to rows (3-6). Shift the existing data down first,
# working from the bottom up so we don't clobber rows before copying them.
for ($r = 5; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 5).Value  = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($dest, 14).Value = $ws.Cells.Item($r, 14).Value()
    $ws.Cells.Item($dest, 16).Value = $ws.Cells.Item($r, 16).Value()
    $ws.Cells.Item($dest, 20).Value = $ws.Cells.Item($r, 20).Value()
}

# Populate the new row 2 with the new transaction details.
$ws.Cells.Item(2, 5).Value  = "Withdrawal"
$ws.Cells.Item(2, 14).Value = "Crypto"
$ws.Cells.Item(2, 16).Value = "ETH"
$ws.Cells.Item(2, 20).Value = 999.98659999999995

# Match the final selection left behind in the workbook.
$ws.Range("Q2:S6").Select()
